# "removing 3 buggy pca methods"
# Net OOXML effect: append one new slide (slide3.xml) at the end of the
# deck. The surviving slide contains a single, empty "Content
# Placeholder" (idx=1, no type) whose shape id/name land on 5 /
# "Content Placeholder 4" -- the numbering left behind once the title
# and the extra (now-removed) placeholders that used to carry the 3
# buggy PCA methods are gone.

$p = $ppt.ActivePresentation

# 1) Add a brand new, blank slide at the end of the deck.
$newIndex = $p.Slides.Count + 1
$s = $p.Slides.Add($newIndex, 12)   # ppLayoutBlank

# 2) Burn through two throw-away shape ids (2, 3) so that once we graft
#    the "Title and Content" layout's placeholders onto this slide they
#    land on ids 4 (Title) and 5 (Content Placeholder) -- matching the
#    ids left over on the real slide after its other placeholders were
#    deleted.
$burn1 = $s.Shapes.AddTextbox(1, 0, 0, 1, 1)
$burn2 = $s.Shapes.AddTextbox(1, 0, 0, 1, 1)
$burn1.Delete()
$burn2.Delete()

# 3) Apply the "Title and Content" layout so the slide gets a real
#    Title + Content placeholder pair (ids 4 and 5).
$s.CustomLayout = $p.SlideMaster.CustomLayouts.Item(2)

# 4) Drop the Title placeholder -- only the content placeholder survives.
$s.Shapes.Item(1).Delete()

# 5) Rename the surviving placeholder to match its id (5 -> "... 4").
$content = $s.Shapes.Item(1)
$content.Name = "Content Placeholder 4"
